# Auto-applies the meteocat daily-summary refresh for 2026-02-07 02:49 UTC run.
# For each changed cell: refreshed DATA_EXTRACCIO timestamp (col E) plus updated
# observation values (HUMITAT_MITJANA_DIA, PRESSIO_ATMOSFERICA, RATXA_VENT_MAX,
# TEMPERATURA_* columns) as scraped from meteo.cat.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "'2026-02-07 02:47:50"
$ws.Range("H2").Value = "'97%"

# Row 3
$ws.Range("E3").Value = "'2026-02-07 02:47:52"
$ws.Range("H3").Value = "'96%"
$ws.Range("N3").Value = "-6.4 °C 2:29 TU"
$ws.Range("O3").Value = "-4.9 °C"

# Row 4
$ws.Range("E4").Value = "'2026-02-07 02:47:55"
$ws.Range("H4").Value = "'55%"
$ws.Range("O4").Value = "12.0 °C"

# Row 5
$ws.Range("E5").Value = "'2026-02-07 02:47:57"
$ws.Range("N5").Value = "7.8 °C 2:24 TU"
$ws.Range("O5").Value = "9.4 °C"

# Row 6
$ws.Range("E6").Value = "'2026-02-07 02:48:00"
$ws.Range("H6").Value = "'58%"
$ws.Range("J6").Value = "1002.1 hPa"
$ws.Range("L6").Value = "27.0 km/h - 260º 2:29 TU"
$ws.Range("N6").Value = "11.4 °C 2:00 TU"
$ws.Range("O6").Value = "12.2 °C"

# Row 7
$ws.Range("E7").Value = "'2026-02-07 02:48:02"
$ws.Range("H7").Value = "'73%"
$ws.Range("J7").Value = "1002.0 hPa"
$ws.Range("N7").Value = "7.4 °C 2:22 TU"
$ws.Range("O7").Value = "8.2 °C"

# Row 8
$ws.Range("E8").Value = "'2026-02-07 02:48:05"
$ws.Range("L8").Value = "7.2 km/h - 78º 2:18 TU"
$ws.Range("N8").Value = "4.2 °C 2:14 TU"
$ws.Range("O8").Value = "5.1 °C"

# Row 9
$ws.Range("E9").Value = "'2026-02-07 02:48:07"
$ws.Range("N9").Value = "2.0 °C 2:25 TU"
$ws.Range("O9").Value = "2.8 °C"

# Row 10
$ws.Range("E10").Value = "'2026-02-07 02:48:10"
$ws.Range("O10").Value = "7.2 °C"

# Row 11
$ws.Range("E11").Value = "'2026-02-07 02:48:12"
$ws.Range("N11").Value = "0.7 °C 2:29 TU"
$ws.Range("O11").Value = "1.4 °C"

# Row 12
$ws.Range("E12").Value = "'2026-02-07 02:48:14"
$ws.Range("N12").Value = "8.8 °C 2:29 TU"
$ws.Range("O12").Value = "10.5 °C"

# Row 13
$ws.Range("E13").Value = "'2026-02-07 02:48:17"
$ws.Range("H13").Value = "'90%"
$ws.Range("O13").Value = "7.2 °C"

# Row 14
$ws.Range("E14").Value = "'2026-02-07 02:48:19"
$ws.Range("H14").Value = "'86%"

# Row 15
$ws.Range("E15").Value = "'2026-02-07 02:48:21"
$ws.Range("J15").Value = "1001.0 hPa"
$ws.Range("N15").Value = "6.4 °C 2:26 TU"
$ws.Range("O15").Value = "8.2 °C"

# Row 16
$ws.Range("E16").Value = "'2026-02-07 02:48:24"
$ws.Range("H16").Value = "'86%"
$ws.Range("N16").Value = "2.7 °C 2:21 TU"
$ws.Range("O16").Value = "3.8 °C"

# Row 17
$ws.Range("E17").Value = "'2026-02-07 02:48:26"

# Row 18
$ws.Range("E18").Value = "'2026-02-07 02:48:29"
$ws.Range("N18").Value = "-7.2 °C 2:29 TU"
$ws.Range("O18").Value = "-6.3 °C"

# Row 19
$ws.Range("E19").Value = "'2026-02-07 02:48:31"
$ws.Range("J19").Value = "1005.3 hPa"
$ws.Range("L19").Value = "13.3 km/h - 231º 2:25 TU"
$ws.Range("N19").Value = "4.1 °C 2:21 TU"
$ws.Range("O19").Value = "5.1 °C"

# Row 20
$ws.Range("E20").Value = "'2026-02-07 02:48:34"
$ws.Range("H20").Value = "'90%"
$ws.Range("L20").Value = "18.0 km/h - 258º 2:28 TU"
$ws.Range("M20").Value = "-3.7 °C 2:27 TU"
$ws.Range("N20").Value = "-4.8 °C 2:08 TU"

# Row 21
$ws.Range("E21").Value = "'2026-02-07 02:48:36"
$ws.Range("H21").Value = "'63%"
$ws.Range("J21").Value = "1000.7 hPa"
$ws.Range("N21").Value = "6.5 °C 2:29 TU"
$ws.Range("O21").Value = "9.3 °C"

# Row 22
$ws.Range("E22").Value = "'2026-02-07 02:48:39"
$ws.Range("H22").Value = "'95%"
$ws.Range("N22").Value = "5.0 °C 2:22 TU"
$ws.Range("O22").Value = "6.0 °C"

# Row 23
$ws.Range("E23").Value = "'2026-02-07 02:48:41"
$ws.Range("M23").Value = "8.2 °C 2:18 TU"
$ws.Range("N23").Value = "7.6 °C 2:29 TU"

# Row 24
$ws.Range("E24").Value = "'2026-02-07 02:48:44"
$ws.Range("L24").Value = "17.6 km/h - 345º 2:29 TU"
$ws.Range("N24").Value = "9.9 °C 2:19 TU"

# Row 25
$ws.Range("E25").Value = "'2026-02-07 02:48:46"
$ws.Range("J25").Value = "1004.8 hPa"
$ws.Range("N25").Value = "0.3 °C 2:19 TU"
$ws.Range("O25").Value = "0.8 °C"

# Row 26
$ws.Range("E26").Value = "'2026-02-07 02:48:49"
$ws.Range("H26").Value = "'76%"

# Row 27
$ws.Range("E27").Value = "'2026-02-07 02:48:51"
$ws.Range("O27").Value = "8.4 °C"

# Row 28
$ws.Range("E28").Value = "'2026-02-07 02:48:54"
$ws.Range("H28").Value = "'84%"
$ws.Range("J28").Value = "1002.9 hPa"
$ws.Range("N28").Value = "3.2 °C 2:29 TU"
$ws.Range("O28").Value = "4.2 °C"

# Row 29
$ws.Range("E29").Value = "'2026-02-07 02:48:56"
$ws.Range("H29").Value = "'55%"
$ws.Range("N29").Value = "10.4 °C 2:03 TU"
$ws.Range("O29").Value = "11.9 °C"

# Row 30
$ws.Range("E30").Value = "'2026-02-07 02:48:58"
$ws.Range("H30").Value = "'81%"
$ws.Range("L30").Value = "31.7 km/h - 338º 2:11 TU"
$ws.Range("N30").Value = "-5.7 °C 2:24 TU"
$ws.Range("O30").Value = "-4.5 °C"

# Row 31
$ws.Range("E31").Value = "'2026-02-07 02:49:01"
$ws.Range("N31").Value = "3.5 °C 2:29 TU"

# Row 32
$ws.Range("E32").Value = "'2026-02-07 02:49:03"
$ws.Range("J32").Value = "1003.6 hPa"

# Row 33
$ws.Range("E33").Value = "'2026-02-07 02:49:06"
$ws.Range("H33").Value = "'89%"
$ws.Range("O33").Value = "7.8 °C"

# Row 34
$ws.Range("E34").Value = "'2026-02-07 02:49:08"
$ws.Range("L34").Value = "23.4 km/h - 231º 2:13 TU"
$ws.Range("N34").Value = "6.1 °C 2:29 TU"
$ws.Range("O34").Value = "7.2 °C"

# Row 35
$ws.Range("E35").Value = "'2026-02-07 02:49:11"
$ws.Range("N35").Value = "-4.3 °C 2:29 TU"

# Row 36
$ws.Range("E36").Value = "'2026-02-07 02:49:13"
$ws.Range("J36").Value = "1005.8 hPa"
$ws.Range("N36").Value = "4.5 °C 2:25 TU"
$ws.Range("O36").Value = "4.8 °C"
